$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 81. This shifts the existing rows 81:202
# down to 82:203, carrying their values/formatting with them, and leaves
# row 81 blank (inheriting formatting from the row above).
$ws.Rows("81").Insert()

# Fill in the new record for row 81 (a weekly "Perejil" price entry for
# Vega Modelo de Temuco, Provincia de Cautin).
$ws.Range("A81").Value = 10
$ws.Range("B81").Value = "Vega Modelo de Temuco"
$ws.Range("C81").Value = "La Araucanía"
$ws.Range("D81").Value = 44477
$ws.Range("E81").Value = 9
$ws.Range("F81").Value = 100112044
$ws.Range("G81").Value = "Perejil"
$ws.Range("H81").Value = "Sin especificar"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 30
$ws.Range("K81").Value = 4000
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = 4000
$ws.Range("N81").Value = "$/docena de atados (3 kilos)"
$ws.Range("O81").Value = "Provincia de Cautín"
$ws.Range("P81").Value = 1333
$ws.Range("Q81").Value = 3
$ws.Range("R81").Value = "Hortaliza"
